# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 1673
    "F4"  = 770
    "F5"  = 1115
    "F6"  = 29
    "F7"  = 11776
    "F8"  = 36
    "F9"  = 96
    "F11" = 394
    "F12" = 1104
    "F13" = 837
    "F14" = 13430
    "F15" = 13311
    "F16" = 37
    "F21" = 90
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
